$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.061.38'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '1.642.91'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  -0.51%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '218.06'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.508'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.01%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.54%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.255'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.69%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.0625'
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '20.30'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +5.36%  '
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").Value = '1.871.27'
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = '1.640.63'
$ws.Range("E13").Value = '  -0.08%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.12'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.538'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.90%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '67.36'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.12%  '
$ws.Range("D17").Value = '27.036.13'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  +0.73%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '221.02'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.44%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +3.23%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.43'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.34%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.11%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.22'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.09%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '147.46'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.36%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.39'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.71%  '
$ws.Range("E28").Value = '  +0.99%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.82'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.41%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0506'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.18'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("E32").Value = '  -0.53%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.87%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.58'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '1.272.31'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("E36").Value = '  -0.03%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0178'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.41%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.546'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.81%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.844'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.96%  '
$ws.Range("E40").Value = '  -0.48%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.809'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.59%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.37'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("D43").Value = '1.781.46'
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '62.89'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.43%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '92.68'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("E48").Value = '  +0.64%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0513'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("E50").Value = '  +1.49%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0971'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.16%  '
